$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.188.04'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.599.49'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '540.20'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.81%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.02'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.565'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '3.062.13'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Value = '59.105.05'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.592.69'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000133'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '342.05'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.58'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '0.0₃0736'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.69'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +8.12%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.80'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.61%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '149.38'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  -0.68%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '37.11'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.30%  '
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('E40').Value = '  +1.20%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '273.76'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.74'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('D48').Value = '1.937.45'
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.43'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '111.32'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.21%  '
